# Updated cryptos list - refresh Price (col D) and Volume(1h) (col E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value,
        [bool]$ForceText = $false
    )
    $cell = $ws.Range($Address)
    if ($ForceText) {
        # Ensure values that look numeric (e.g. "567.97") are stored as
        # literal text, matching the original inline-string cell content
        # instead of being auto-converted into a floating point number.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $Value
}

# --- Column D (Price) updates ---
Set-TextValue "D2"  "61.092.07"
Set-TextValue "D3"  "2.401.20"
Set-TextValue "D5"  "567.97"    $true
Set-TextValue "D6"  "141.94"    $true
Set-TextValue "D9"  "2.409.53"
Set-TextValue "D12" "5.20"      $true
Set-TextValue "D14" "26.46"     $true
Set-TextValue "D16" "2.840.12"
Set-TextValue "D17" "60.976.15"
Set-TextValue "D18" "2.413.01"
Set-TextValue "D19" "8.06"      $true
Set-TextValue "D21" "323.93"    $true
Set-TextValue "D26" "65.16"     $true
Set-TextValue "D27" "589.20"    $true
Set-TextValue "D28" "8.22"      $true
Set-TextValue "D31" "8.01"      $true
Set-TextValue "D37" "153.24"    $true
Set-TextValue "D38" "0.373"     $true
Set-TextValue "D40" "18.35"     $true
Set-TextValue "D44" "41.73"     $true
Set-TextValue "D46" "0.0₆0283"
Set-TextValue "D47" "141.60"    $true
Set-TextValue "D48" "3.53"      $true
Set-TextValue "D50" "19.69"     $true
Set-TextValue "D51" "0.0511"    $true

# --- Column E (Volume 1h) updates ---
Set-TextValue "E2"  "  +0.01%  "
Set-TextValue "E4"  "  +0.46%  "
Set-TextValue "E5"  "  -0.17%  "
Set-TextValue "E6"  "  +1.78%  "
Set-TextValue "E7"  "  -0.37%  "
Set-TextValue "E8"  "  +2.12%  "
Set-TextValue "E9"  "  +0.16%  "
Set-TextValue "E11" "  -0.29%  "
Set-TextValue "E12" "  +2.60%  "
Set-TextValue "E13" "  +2.46%  "
Set-TextValue "E14" "  +1.57%  "
Set-TextValue "E15" "  -0.19%  "
Set-TextValue "E16" "  -0.51%  "
Set-TextValue "E17" "  -0.05%  "
Set-TextValue "E18" "  +0.65%  "
Set-TextValue "E19" "  +2.58%  "
Set-TextValue "E20" "  +0.49%  "
Set-TextValue "E21" "  +0.34%  "
Set-TextValue "E22" "  +1.06%  "
Set-TextValue "E23" "  -0.22%  "
Set-TextValue "E24" "  -0.25%  "
Set-TextValue "E25" "  +4.58%  "
Set-TextValue "E26" "  +0.80%  "
Set-TextValue "E27" "  +1.17%  "
Set-TextValue "E28" "  -0.27%  "
Set-TextValue "E29" "  +1.92%  "
Set-TextValue "E31" "  +1.64%  "
Set-TextValue "E32" "  +1.94%  "
Set-TextValue "E33" "  -0.55%  "
Set-TextValue "E34" "  +0.18%  "
Set-TextValue "E35" "  +3.98%  "
Set-TextValue "E36" "  -0.65%  "
Set-TextValue "E37" "  +0.81%  "
Set-TextValue "E38" "  +1.22%  "
Set-TextValue "E39" "  +0.26%  "
Set-TextValue "E41" "  +1.91%  "
Set-TextValue "E43" "  +1.23%  "
Set-TextValue "E44" "  +1.47%  "
Set-TextValue "E45" "  +6.87%  "
Set-TextValue "E46" "  +2.52%  "
Set-TextValue "E47" "  -1.07%  "
Set-TextValue "E48" "  +0.59%  "
Set-TextValue "E49" "  +0.99%  "
Set-TextValue "E50" "  +1.04%  "
Set-TextValue "E51" "  +1.80%  "

Write-Host "Updated cryptos list values."
